# Auto update Excel log
# Appends 10 new PRESENCE_DETECTED rows (38-47) to the "mmWave" sheet,
# matching the sensor-log pattern already present in rows 2-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @{ Row = 38; Time = "15:57:51" },
    @{ Row = 39; Time = "15:58:02" },
    @{ Row = 40; Time = "15:58:12" },
    @{ Row = 41; Time = "15:58:23" },
    @{ Row = 42; Time = "15:58:44" },
    @{ Row = 43; Time = "15:58:54" },
    @{ Row = 44; Time = "15:59:05" },
    @{ Row = 45; Time = "15:59:15" },
    @{ Row = 46; Time = "15:59:25" },
    @{ Row = 47; Time = "15:59:36" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Column A holds a plain date-looking string ("2026-02-01"). Force the
    # cell to Text format first so Excel keeps it as literal text instead of
    # silently converting it to a date serial number, then restore the
    # default "Normal" style so no stray number-format is left on the cell.
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = "2026-02-01"
    $ws.Range("A$r").Style = "Normal"

    $ws.Range("B$r").Value = $entry.Time
    $ws.Range("C$r").Value = "15:00"
    $ws.Range("D$r").Value = "Living Room"
    $ws.Range("E$r").Value = "PRESENCE_DETECTED"
    $ws.Range("F$r").Value = "Active"
}
